# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "K" (strikeouts) values replacing the old "Strike#" values in column G,
# for data rows 2 through 40 (row 1 is the header row).
$kValues = @(9, 9, 5, 8, 6, 7, 5, 8, 2, 4, 4, 1, 2, 1, 0, 0, 1, 3, 4, 2, 1, 2, 3, 0, 1, 3, 1, 0, 3, 3, 1, 0, 3, 1, 0, 1, 1, 2, 3)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
